$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster data (header + 17 players), replacing the previous 18-row list.
$data = @(
    @("Oyuncu Adı", "Pozisyon", "Takım"),
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Kentavious Caldwell-Pope", "SG,SF", "Orlando Magic"),
    @("Harrison Barnes", "SF,PF", "San Antonio Spurs"),
    @("Malik Monk", "SG,SF", "Sacramento Kings"),
    @("Julius Randle", "PF", "Minnesota Timberwolves"),
    @("Anthony Davis", "PF,C", "Los Angeles Lakers"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Nick Richards", "C", "Charlotte Hornets"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Quentin Grimes", "SG,SF", "Dallas Mavericks"),
    @("Brandon Miller", "SG,SF", "Charlotte Hornets"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets"),
    @("Malcolm Brogdon", "PG,SG", "Washington Wizards")
)

# Clear old contents (not formatting) since the new table is one row shorter.
$ws.Range("A1:C19").ClearContents()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
